# The deck ships two theme parts: ppt/theme/theme1.xml (the main/deck theme,
# wired to the slide master -> all slides) and ppt/theme/theme2.xml (wired to
# the notes master only). The commit swaps their contents: theme1.xml becomes
# the stock "Office Theme" colour scheme, theme2.xml becomes the "Integral"
# colour scheme that theme1.xml used to hold. The font scheme and format
# scheme (fills/lines/effects) are byte-identical between the two themes, so
# the only observable change is the 12-slot theme colour scheme (and the
# cosmetic theme/clrScheme "name" attributes, which PowerPoint's object model
# does not expose for editing).
#
# Apply the reachable half of that swap: recolor the deck's theme (the one
# backing the slide master / every slide) from "Integral" to the "Office"
# palette.

function Get-RGBFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Order of ThemeColorScheme.Item(n) is: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = Get-RGBFromHex $officeThemeColors[$i - 1]
}
